# Fixing product embedding update:
# - Row 2 ("Ayam Bakar") had the wrong description (a drink description
#   left over from copy/paste). Replace it with the correct description
#   for grilled chicken.
# - Update the active selection on the Products sheet to D8 (where the
#   fix was reviewed) instead of the old G5 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

$ws.Range("D2").Value = "Makanan bergizi terbuat dari ayam yang dibakar"

$ws.Activate()
$ws.Range("D8").Select()
